$d = $word.ActiveDocument

$pairs = @(
    @("402×6=", "858×3="),
    @("884×5=", "374×2="),
    @("510×2=", "956×3="),
    @("177×6=", "418×9="),
    @("849×2=", "214×9="),
    @("765×9=", "806×8="),
    @("901×6=", "847×7="),
    @("850×7=", "953×3="),
    @("675×9=", "922×4="),
    @("291×5=", "502×8="),
    @("269×4=", "573×7="),
    @("414×8=", "758×6="),
    @("838×5=", "414×5="),
    @("341×6=", "616×4="),
    @("696×8=", "693×5="),
    @("680×4=", "284×4="),
    @("837×2=", "853×2="),
    @("431×4=", "233×3="),
    @("180×3=", "468×7="),
    @("548×3=", "847×7="),
    @("667×6=", "426×2="),
    @("142×3=", "208×4="),
    @("192×3=", "514×2="),
    @("142×4=", "838×3="),
    @("111×2=", "270×2=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
